$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.656.25'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +3.54%  '
$c.ClearFormats()

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.873.67'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +3.22%  '
$c.ClearFormats()

# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.9988'
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.25%  '
$c.ClearFormats()

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '282.33'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +2.03%  '
$c.ClearFormats()

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.9985'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.ClearFormats()

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5164'
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +3.43%  '
$c.ClearFormats()

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3542'
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +1.17%  '
$c.ClearFormats()

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '45.31'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +2.49%  '
$c.ClearFormats()

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.07035'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +6.14%  '
$c.ClearFormats()

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '20.15'
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.56%  '
$c.ClearFormats()

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.8199'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -2.82%  '
$c.ClearFormats()

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.07749'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '
$c.ClearFormats()

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.875.49'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +3.37%  '
$c.ClearFormats()

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.144'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +1.95%  '
$c.ClearFormats()

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '89.58'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +2.40%  '
$c.ClearFormats()

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.9986'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c.ClearFormats()

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '14.44'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +3.94%  '
$c.ClearFormats()

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.000008156'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +2.13%  '
$c.ClearFormats()

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.9985'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -0.27%  '
$c.ClearFormats()

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '26.688.11'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +3.38%  '
$c.ClearFormats()

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.809'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +1.60%  '
$c.ClearFormats()

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.14'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +1.42%  '
$c.ClearFormats()

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.250'
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +2.91%  '
$c.ClearFormats()

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.436'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +15.86%  '
$c.ClearFormats()

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '145.33'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +3.19%  '
$c.ClearFormats()

# Row 27
$c = $ws.Range('B27')
$c.NumberFormat = '@'
$c.Value = 'EthereumClassic'
$c.ClearFormats()
$c = $ws.Range('C27')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '17.39'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +3.21%  '
$c.ClearFormats()

# Row 28
$c = $ws.Range('B28')
$c.NumberFormat = '@'
$c.Value = 'Toncoin'
$c.ClearFormats()
$c = $ws.Range('C28')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.661'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -0.31%  '
$c.ClearFormats()

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '110.93'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +1.95%  '
$c.ClearFormats()

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.412'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +1.90%  '
$c.ClearFormats()

# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.362'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +3.85%  '
$c.ClearFormats()

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.08833'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +0.85%  '
$c.ClearFormats()

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.04920'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +1.95%  '
$c.ClearFormats()

# Row 34
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +4.51%  '
$c.ClearFormats()

# Row 35
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.7488'
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +1.44%  '
$c.ClearFormats()

# Row 36
$c = $ws.Range('B36')
$c.NumberFormat = '@'
$c.Value = 'MXToken'
$c.ClearFormats()
$c = $ws.Range('C36')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.295'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +7.39%  '
$c.ClearFormats()

# Row 37
$c = $ws.Range('B37')
$c.NumberFormat = '@'
$c.Value = 'HuobiToken'
$c.ClearFormats()
$c = $ws.Range('C37')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.ClearFormats()
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.862'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -0.48%  '
$c.ClearFormats()

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.431'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -1.86%  '
$c.ClearFormats()

# Row 39
$c = $ws.Range('B39')
$c.NumberFormat = '@'
$c.Value = 'TheSandbox'
$c.ClearFormats()
$c = $ws.Range('C39')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.5287'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -0.15%  '
$c.ClearFormats()

# Row 40
$c = $ws.Range('B40')
$c.NumberFormat = '@'
$c.Value = 'VeChain'
$c.ClearFormats()
$c = $ws.Range('C40')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.ClearFormats()
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.01878'
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +0.94%  '
$c.ClearFormats()

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.9718'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.ClearFormats()

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '116.98'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +4.23%  '
$c.ClearFormats()

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '6.318'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +1.53%  '
$c.ClearFormats()

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '8.200'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +0.84%  '
$c.ClearFormats()

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.9983'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.32%  '
$c.ClearFormats()

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4609'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -1.72%  '
$c.ClearFormats()

# Row 47
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.1365'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -1.51%  '
$c.ClearFormats()

# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.514'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +3.41%  '
$c.ClearFormats()

# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '36.64'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +2.98%  '
$c.ClearFormats()

# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.518'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +2.11%  '
$c.ClearFormats()

# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.05926'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +0.83%  '
$c.ClearFormats()
